$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.613.52"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "3.780.48"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'595.20"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'166.36"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "3.779.83"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "'6.35"
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "'36.25"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "4.415.96"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "3.778.53"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "'18.43"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "67.600.26"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -5.59%  "
$ws.Range("D22").Value = "'457.07"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'0.0000155"
$ws.Range("E24").Value = "  +8.00%  "
$ws.Range("D25").Value = "'83.35"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").Value = "'11.92"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'2.78"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'29.80"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "3.734.11"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "'0.137"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").Value = "'0.994"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D44").Value = "'45.29"
$ws.Range("E44").Value = "  +5.84%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "'47.09"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").Value = "'148.27"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("D50").Value = "'389.50"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "'25.62"
$ws.Range("E51").Value = "  +0.75%  "
